# Applies six copy-edit corrections to the Partner Playbook template,
# using Word's Find/Replace (wdReplaceAll) via the Find.Execute COM call.
# Signature reminder:
#   Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#                MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#                Format, ReplaceWith, Replace)

$d = $word.ActiveDocument
$emDash = [char]0x2014
$wdFindContinue = 1
$wdReplaceAll = 2

$edits = @(
    @{
        Old = "All partners collaborate to define the roles needed, responsibilities, and how to ensure their team members are supported."
        New = "All partners collaborate to define the proposal team roles, responsibilities, and how to ensure their team members are supported."
    },
    @{
        Old = "Prime-level partners collaborate and participate across all phases of bid preparation $emDash scheduling early with special care and attention to reduce effort outside of business hours."
        New = "Prime-level partners collaborate and participate across all phases of bid preparation $emDash scheduling early with special care to reduce effort outside of business hours."
    },
    @{
        Old = "All partners discuss rates, labor category descriptions and discount strategy options early to ensure we are aligned on how staffing impacts pricing."
        New = "All partners discuss rates, labor category descriptions and discount strategies early to ensure the team is aligned on how staffing impacts pricing."
    },
    @{
        Old = "All partners hold each other accountable by openly communicating status, progress, and blockers. "
        New = "All partners hold each other accountable by openly communicating status, intel, progress, and blockers. "
    },
    @{
        Old = "Prime-level partners communicate with client stakeholders (e.g., Contracting Officer), and seek opportunities for sub engagement with stakeholders. "
        New = "Prime-level partners communicate with client stakeholders (e.g., Contracting Officer), and seek opportunities for subcontractors to engage with stakeholders. "
    },
    @{
        Old = "Partners collaborate on co-marketing and strategic positioning to support cost modification and recompete success. "
        New = "Partners collaborate on co-marketing and strategic positioning to support growth opportunities (e.g., cost modification and recompete success). "
    }
)

foreach ($edit in $edits) {
    $ok = $d.Content.Find.Execute(
        $edit.Old, $true, $false, $false, $false, $false, $true,
        $wdFindContinue, $false, $edit.New, $wdReplaceAll)
    Write-Output "replaced='$ok' old='$($edit.Old.Substring(0, [Math]::Min(50, $edit.Old.Length)))...'"
}
